$d = $word.ActiveDocument

$replacements = @(
    @("66×72=", "20×39="),
    @("22×90=", "61×28="),
    @("89×36=", "18×60="),
    @("72×96=", "98×68="),
    @("79×69=", "64×73="),
    @("17×74=", "75×37="),
    @("19×37=", "48×63="),
    @("13×35=", "74×69="),
    @("41×23=", "41×16="),
    @("62×77=", "97×90="),
    @("45×23=", "98×49="),
    @("15×44=", "85×40="),
    @("13×93=", "71×89="),
    @("80×18=", "42×31="),
    @("60×59=", "85×99="),
    @("64×34=", "46×77="),
    @("73×81=", "85×61="),
    @("24×20=", "82×47="),
    @("69×98=", "56×64="),
    @("36×64=", "14×31="),
    @("51×34=", "56×69="),
    @("24×40=", "47×51="),
    @("64×79=", "70×23="),
    @("53×56=", "36×52="),
    @("18×66=", "68×38=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
